$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit: a value was typed into cell A1.
$ws.Range("A1").Value = "change added"

# The sheet's print setup was also saved as A4 portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
